$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing AgTests (F) / AgPosit (G) values for rows 279-428 ---
$ws.Range("F279").Value = 42176
$ws.Range("G279").Value = 2984
$ws.Range("F336").Value = 81542
$ws.Range("F338").Value = 221029
$ws.Range("G338").Value = 3044
$ws.Range("F339").Value = 662655
$ws.Range("F341").Value = 283207
$ws.Range("G341").Value = 3608
$ws.Range("F342").Value = 178532
$ws.Range("F348").Value = 232768
$ws.Range("F349").Value = 159359
$ws.Range("F350").Value = 127002
$ws.Range("F351").Value = 150474
$ws.Range("G351").Value = 2806
$ws.Range("F352").Value = 307353
$ws.Range("G352").Value = 3547
$ws.Range("F353").Value = 723546
$ws.Range("G353").Value = 5290
$ws.Range("F355").Value = 221931
$ws.Range("G355").Value = 3434
$ws.Range("F356").Value = 159949
$ws.Range("F357").Value = 138195
$ws.Range("G357").Value = 3012
$ws.Range("F358").Value = 158767
$ws.Range("G358").Value = 2610
$ws.Range("F359").Value = 321150
$ws.Range("G359").Value = 3334
$ws.Range("F360").Value = 749709
$ws.Range("F365").Value = 184591
$ws.Range("F371").Value = 160004
$ws.Range("F372").Value = 178378
$ws.Range("G372").Value = 1852
$ws.Range("F374").Value = 773531
$ws.Range("G374").Value = 3419
$ws.Range("F377").Value = 176537
$ws.Range("G377").Value = 1823
$ws.Range("F380").Value = 344416
$ws.Range("G380").Value = 2022
$ws.Range("F381").Value = 745647
$ws.Range("F383").Value = 220778
$ws.Range("F385").Value = 150744
$ws.Range("G385").Value = 1408
$ws.Range("F387").Value = 351552
$ws.Range("F391").Value = 177241
$ws.Range("F392").Value = 220865
$ws.Range("G392").Value = 1216
$ws.Range("F393").Value = 308126
$ws.Range("F394").Value = 166172
$ws.Range("F395").Value = 750792
$ws.Range("F396").Value = 164662
$ws.Range("F397").Value = 108022
$ws.Range("F398").Value = 298177
$ws.Range("G398").Value = 1467
$ws.Range("F399").Value = 200436
$ws.Range("G399").Value = 966
$ws.Range("F400").Value = 150762
$ws.Range("G400").Value = 767
$ws.Range("F401").Value = 273355
$ws.Range("F402").Value = 716948
$ws.Range("F403").Value = 351927
$ws.Range("G403").Value = 730
$ws.Range("F404").Value = 224905
$ws.Range("G404").Value = 912
$ws.Range("F405").Value = 173903
$ws.Range("F406").Value = 170839
$ws.Range("F408").Value = 303656
$ws.Range("F409").Value = 703004
$ws.Range("F410").Value = 363434
$ws.Range("F411").Value = 225001
$ws.Range("F412").Value = 175804
$ws.Range("G412").Value = 645
$ws.Range("F413").Value = 148919
$ws.Range("F414").Value = 146763
$ws.Range("G414").Value = 557
$ws.Range("F415").Value = 304867
$ws.Range("G415").Value = 693
$ws.Range("F416").Value = 658759
$ws.Range("F417").Value = 332534
$ws.Range("F418").Value = 200363
$ws.Range("F419").Value = 147463
$ws.Range("G419").Value = 500
$ws.Range("F421").Value = 150371
$ws.Range("G421").Value = 525
$ws.Range("F422").Value = 293474
$ws.Range("F423").Value = 432090
$ws.Range("F424").Value = 255004
$ws.Range("G424").Value = 485
$ws.Range("F425").Value = 136504
$ws.Range("G425").Value = 542
$ws.Range("F426").Value = 104543
$ws.Range("G426").Value = 384
$ws.Range("F427").Value = 89031
$ws.Range("G427").Value = 358
$ws.Range("F428").Value = 99209
$ws.Range("G428").Value = 374

# --- Append two new rows for 2021-05-07 and 2021-05-08 ---
$ws.Cells.Item(429, 1).Value = 44323
$ws.Cells.Item(429, 2).Value = 385395
$ws.Cells.Item(429, 3).Value = 12132
$ws.Cells.Item(429, 4).Value = 373
$ws.Cells.Item(429, 5).Value = 11990
$ws.Cells.Item(429, 6).Value = 160230
$ws.Cells.Item(429, 7).Value = 585

$ws.Cells.Item(430, 1).Value = 44324
$ws.Cells.Item(430, 2).Value = 385475
$ws.Cells.Item(430, 3).Value = 3666
$ws.Cells.Item(430, 4).Value = 80
$ws.Cells.Item(430, 5).Value = 12019
$ws.Cells.Item(430, 6).Value = 152952
$ws.Cells.Item(430, 7).Value = 249
